$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 188
$ws.Range("A188").Value = 111973724
$ws.Range("B188").Value = 88489
$ws.Range("C188").Value = "Ovaliderad"
$ws.Range("D188").Value = "NT"
$ws.Range("E188").Value = 1962
$ws.Range("F188").Value = "Vaddporing"
$ws.Range("G188").Value = "Anomoporia kamtschatica"
$ws.Range("H188").Value = "(Parmasto) Bondartseva"
$ws.Range("I188").Value = ""
$ws.Range("J188").ClearContents()
$ws.Range("P188").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q188").Value = 437892.4817196695
$ws.Range("R188").Value = 6953090.40240525
$ws.Range("S188").Value = 10
$ws.Range("T188").Value = "Jämtland"
$ws.Range("U188").Value = "Härjedalen"
$ws.Range("V188").Value = "Jämtland"
$ws.Range("W188").Value = "Vemdalen"
$ws.Range("Y188").Value = "2023-09-07"
$ws.Range("Z188").Value = "00:00"
$ws.Range("AA188").Value = "2023-09-07"
$ws.Range("AB188").Value = "00:00"
$ws.Range("AD188").Value = $false
$ws.Range("AE188").Value = $false
$ws.Range("AG188").Value = $false
$ws.Range("AI188").Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"
$ws.Range("AN188").Value = 1
$ws.Range("AO188").Value = "1 substratenheter # under gammal tallåga"
$ws.Range("AT188").Value = ""
$ws.Range("AW188").Value = "Magnus Andersson"
$ws.Range("AX188").Value = "Magnus Andersson"
$ws.Range("AY188").Value = "SCA Skog Naturvärdesinventering"

# Row 189
$ws.Range("A189").Value = 111973651
$ws.Range("B189").Value = 90682
$ws.Range("C189").Value = "Ovaliderad"
$ws.Range("D189").Value = "NT"
$ws.Range("E189").Value = 2059
$ws.Range("F189").Value = "Skrovlig taggsvamp"
$ws.Range("G189").Value = "Hydnellum scabrosum"
$ws.Range("H189").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("I189").Value = "1"
$ws.Range("J189").Value = "mycel"
$ws.Range("P189").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q189").Value = 438594.9557070844
$ws.Range("R189").Value = 6953584.041166852
$ws.Range("S189").Value = 10
$ws.Range("T189").Value = "Jämtland"
$ws.Range("U189").Value = "Härjedalen"
$ws.Range("V189").Value = "Jämtland"
$ws.Range("W189").Value = "Vemdalen"
$ws.Range("Y189").Value = "2023-09-07"
$ws.Range("Z189").Value = "00:00"
$ws.Range("AA189").Value = "2023-09-07"
$ws.Range("AB189").Value = "00:00"
$ws.Range("AD189").Value = $false
$ws.Range("AE189").Value = $false
$ws.Range("AG189").Value = $false
$ws.Range("AI189").Value = "äldre fattigristallskog på torr moränmark"
$ws.Range("AN189").ClearContents()
$ws.Range("AO189").ClearContents()
$ws.Range("AT189").Value = ""
$ws.Range("AW189").Value = "Magnus Andersson"
$ws.Range("AX189").Value = "Magnus Andersson"
$ws.Range("AY189").Value = "SCA Skog Naturvärdesinventering"

# Row 190
$ws.Range("A190").Value = 111973764
$ws.Range("B190").Value = 90660
$ws.Range("C190").Value = "Ovaliderad"
$ws.Range("D190").Value = "NT"
$ws.Range("E190").Value = 4362
$ws.Range("F190").Value = "Blå taggsvamp"
$ws.Range("G190").Value = "Hydnellum caeruleum"
$ws.Range("H190").Value = "(Hornem.) P.Karst."
$ws.Range("I190").Value = "1"
$ws.Range("J190").Value = "mycel"
$ws.Range("P190").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q190").Value = 437776.3383109252
$ws.Range("R190").Value = 6953123.809258236
$ws.Range("S190").Value = 10
$ws.Range("T190").Value = "Jämtland"
$ws.Range("U190").Value = "Härjedalen"
$ws.Range("V190").Value = "Jämtland"
$ws.Range("W190").Value = "Vemdalen"
$ws.Range("Y190").Value = "2023-09-07"
$ws.Range("Z190").Value = "00:00"
$ws.Range("AA190").Value = "2023-09-07"
$ws.Range("AB190").Value = "00:00"
$ws.Range("AD190").Value = $false
$ws.Range("AE190").Value = $false
$ws.Range("AG190").Value = $false
$ws.Range("AI190").Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"
$ws.Range("AN190").ClearContents()
$ws.Range("AO190").ClearContents()
$ws.Range("AT190").Value = ""
$ws.Range("AW190").Value = "Magnus Andersson"
$ws.Range("AX190").Value = "Magnus Andersson"
$ws.Range("AY190").Value = "SCA Skog Naturvärdesinventering"

# Row 191
$ws.Range("A191").Value = 111973727
$ws.Range("B191").Value = 90660
$ws.Range("C191").Value = "Ovaliderad"
$ws.Range("D191").Value = "NT"
$ws.Range("E191").Value = 4362
$ws.Range("F191").Value = "Blå taggsvamp"
$ws.Range("G191").Value = "Hydnellum caeruleum"
$ws.Range("H191").Value = "(Hornem.) P.Karst."
$ws.Range("I191").Value = "1"
$ws.Range("J191").Value = "mycel"
$ws.Range("P191").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q191").Value = 437870.6858627723
$ws.Range("R191").Value = 6953178.038412376
$ws.Range("S191").Value = 10
$ws.Range("T191").Value = "Jämtland"
$ws.Range("U191").Value = "Härjedalen"
$ws.Range("V191").Value = "Jämtland"
$ws.Range("W191").Value = "Vemdalen"
$ws.Range("Y191").Value = "2023-09-07"
$ws.Range("Z191").Value = "00:00"
$ws.Range("AA191").Value = "2023-09-07"
$ws.Range("AB191").Value = "00:00"
$ws.Range("AD191").Value = $false
$ws.Range("AE191").Value = $false
$ws.Range("AG191").Value = $false
$ws.Range("AI191").Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"
$ws.Range("AN191").ClearContents()
$ws.Range("AO191").ClearContents()
$ws.Range("AT191").Value = ""
$ws.Range("AW191").Value = "Magnus Andersson"
$ws.Range("AX191").Value = "Magnus Andersson"
$ws.Range("AY191").Value = "SCA Skog Naturvärdesinventering"

# Row 192
$ws.Range("A192").Value = 111973715
$ws.Range("B192").Value = 90678
$ws.Range("C192").Value = "Ovaliderad"
$ws.Range("D192").Value = "LC"
$ws.Range("E192").Value = 4366
$ws.Range("F192").Value = "Skarp dropptaggsvamp"
$ws.Range("G192").Value = "Hydnellum peckii"
$ws.Range("H192").Value = "Banker"
$ws.Range("I192").Value = "1"
$ws.Range("J192").Value = "mycel"
$ws.Range("P192").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q192").Value = 437962.8122493967
$ws.Range("R192").Value = 6953212.111986059
$ws.Range("S192").Value = 10
$ws.Range("T192").Value = "Jämtland"
$ws.Range("U192").Value = "Härjedalen"
$ws.Range("V192").Value = "Jämtland"
$ws.Range("W192").Value = "Vemdalen"
$ws.Range("Y192").Value = "2023-09-07"
$ws.Range("Z192").Value = "00:00"
$ws.Range("AA192").Value = "2023-09-07"
$ws.Range("AB192").Value = "00:00"
$ws.Range("AD192").Value = $false
$ws.Range("AE192").Value = $false
$ws.Range("AG192").Value = $false
$ws.Range("AI192").Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"
$ws.Range("AN192").ClearContents()
$ws.Range("AO192").ClearContents()
$ws.Range("AT192").Value = ""
$ws.Range("AW192").Value = "Magnus Andersson"
$ws.Range("AX192").Value = "Magnus Andersson"
$ws.Range("AY192").Value = "SCA Skog Naturvärdesinventering"

# Row 193
$ws.Range("A193").Value = 111973713
$ws.Range("B193").Value = 90660
$ws.Range("C193").Value = "Ovaliderad"
$ws.Range("D193").Value = "NT"
$ws.Range("E193").Value = 4362
$ws.Range("F193").Value = "Blå taggsvamp"
$ws.Range("G193").Value = "Hydnellum caeruleum"
$ws.Range("H193").Value = "(Hornem.) P.Karst."
$ws.Range("I193").Value = "1"
$ws.Range("J193").Value = "mycel"
$ws.Range("P193").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q193").Value = 438002.4574124058
$ws.Range("R193").Value = 6953193.462733216
$ws.Range("S193").Value = 10
$ws.Range("T193").Value = "Jämtland"
$ws.Range("U193").Value = "Härjedalen"
$ws.Range("V193").Value = "Jämtland"
$ws.Range("W193").Value = "Vemdalen"
$ws.Range("Y193").Value = "2023-09-07"
$ws.Range("Z193").Value = "00:00"
$ws.Range("AA193").Value = "2023-09-07"
$ws.Range("AB193").Value = "00:00"
$ws.Range("AD193").Value = $false
$ws.Range("AE193").Value = $false
$ws.Range("AG193").Value = $false
$ws.Range("AI193").Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"
$ws.Range("AN193").ClearContents()
$ws.Range("AO193").ClearContents()
$ws.Range("AT193").Value = ""
$ws.Range("AW193").Value = "Magnus Andersson"
$ws.Range("AX193").Value = "Magnus Andersson"
$ws.Range("AY193").Value = "SCA Skog Naturvärdesinventering"

# Row 194
$ws.Range("A194").Value = 111973652
$ws.Range("B194").Value = 90660
$ws.Range("C194").Value = "Ovaliderad"
$ws.Range("D194").Value = "NT"
$ws.Range("E194").Value = 4362
$ws.Range("F194").Value = "Blå taggsvamp"
$ws.Range("G194").Value = "Hydnellum caeruleum"
$ws.Range("H194").Value = "(Hornem.) P.Karst."
$ws.Range("I194").Value = "1"
$ws.Range("J194").Value = "mycel"
$ws.Range("P194").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q194").Value = 438537.7865028595
$ws.Range("R194").Value = 6953550.220835418
$ws.Range("S194").Value = 10
$ws.Range("T194").Value = "Jämtland"
$ws.Range("U194").Value = "Härjedalen"
$ws.Range("V194").Value = "Jämtland"
$ws.Range("W194").Value = "Vemdalen"
$ws.Range("Y194").Value = "2023-09-07"
$ws.Range("Z194").Value = "00:00"
$ws.Range("AA194").Value = "2023-09-07"
$ws.Range("AB194").Value = "00:00"
$ws.Range("AD194").Value = $false
$ws.Range("AE194").Value = $false
$ws.Range("AG194").Value = $false
$ws.Range("AI194").Value = "äldre fattigristallskog på torr moränmark"
$ws.Range("AN194").ClearContents()
$ws.Range("AO194").ClearContents()
$ws.Range("AT194").Value = ""
$ws.Range("AW194").Value = "Magnus Andersson"
$ws.Range("AX194").Value = "Magnus Andersson"
$ws.Range("AY194").Value = "SCA Skog Naturvärdesinventering"

# Row 195
$ws.Range("A195").Value = 111973671
$ws.Range("B195").Value = 88032
$ws.Range("C195").Value = "Ovaliderad"
$ws.Range("D195").Value = "VU"
$ws.Range("E195").Value = 6276
$ws.Range("F195").Value = "Goliatmusseron"
$ws.Range("G195").Value = "Tricholoma matsutake"
$ws.Range("H195").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I195").Value = "1"
$ws.Range("J195").Value = "mycel"
$ws.Range("P195").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q195").Value = 438033.4411253001
$ws.Range("R195").Value = 6953252.100307667
$ws.Range("S195").Value = 10
$ws.Range("T195").Value = "Jämtland"
$ws.Range("U195").Value = "Härjedalen"
$ws.Range("V195").Value = "Jämtland"
$ws.Range("W195").Value = "Vemdalen"
$ws.Range("Y195").Value = "2023-09-07"
$ws.Range("Z195").Value = "00:00"
$ws.Range("AA195").Value = "2023-09-07"
$ws.Range("AB195").Value = "00:00"
$ws.Range("AD195").Value = $false
$ws.Range("AE195").Value = $false
$ws.Range("AG195").Value = $false
$ws.Range("AI195").Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"
$ws.Range("AN195").ClearContents()
$ws.Range("AO195").ClearContents()
$ws.Range("AT195").Value = ""
$ws.Range("AW195").Value = "Magnus Andersson"
$ws.Range("AX195").Value = "Magnus Andersson"
$ws.Range("AY195").Value = "SCA Skog Naturvärdesinventering"

# Row 196
$ws.Range("A196").Value = 111973733
$ws.Range("B196").Value = 90660
$ws.Range("C196").Value = "Ovaliderad"
$ws.Range("D196").Value = "NT"
$ws.Range("E196").Value = 4362
$ws.Range("F196").Value = "Blå taggsvamp"
$ws.Range("G196").Value = "Hydnellum caeruleum"
$ws.Range("H196").Value = "(Hornem.) P.Karst."
$ws.Range("I196").Value = "1"
$ws.Range("J196").Value = "mycel"
$ws.Range("P196").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q196").Value = 437876.3203048867
$ws.Range("R196").Value = 6953355.130729643
$ws.Range("S196").Value = 10
$ws.Range("T196").Value = "Jämtland"
$ws.Range("U196").Value = "Härjedalen"
$ws.Range("V196").Value = "Jämtland"
$ws.Range("W196").Value = "Vemdalen"
$ws.Range("Y196").Value = "2023-09-07"
$ws.Range("Z196").Value = "00:00"
$ws.Range("AA196").Value = "2023-09-07"
$ws.Range("AB196").Value = "00:00"
$ws.Range("AD196").Value = $false
$ws.Range("AE196").Value = $false
$ws.Range("AG196").Value = $false
$ws.Range("AI196").Value = "äldre renbetad lingontallskog med lavfläckar på torr moränmark"
$ws.Range("AN196").ClearContents()
$ws.Range("AO196").ClearContents()
$ws.Range("AT196").Value = ""
$ws.Range("AW196").Value = "Magnus Andersson"
$ws.Range("AX196").Value = "Magnus Andersson"
$ws.Range("AY196").Value = "SCA Skog Naturvärdesinventering"

# Row 197
$ws.Range("A197").Value = 111973709
$ws.Range("B197").Value = 90682
$ws.Range("C197").Value = "Ovaliderad"
$ws.Range("D197").Value = "NT"
$ws.Range("E197").Value = 2059
$ws.Range("F197").Value = "Skrovlig taggsvamp"
$ws.Range("G197").Value = "Hydnellum scabrosum"
$ws.Range("H197").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("I197").Value = "1"
$ws.Range("J197").Value = "mycel"
$ws.Range("P197").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q197").Value = 438062.0622929674
$ws.Range("R197").Value = 6953111.545111132
$ws.Range("S197").Value = 10
$ws.Range("T197").Value = "Jämtland"
$ws.Range("U197").Value = "Härjedalen"
$ws.Range("V197").Value = "Jämtland"
$ws.Range("W197").Value = "Vemdalen"
$ws.Range("Y197").Value = "2023-09-07"
$ws.Range("Z197").Value = "00:00"
$ws.Range("AA197").Value = "2023-09-07"
$ws.Range("AB197").Value = "00:00"
$ws.Range("AD197").Value = $false
$ws.Range("AE197").Value = $false
$ws.Range("AG197").Value = $false
$ws.Range("AI197").Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"
$ws.Range("AN197").ClearContents()
$ws.Range("AO197").ClearContents()
$ws.Range("AT197").Value = ""
$ws.Range("AW197").Value = "Magnus Andersson"
$ws.Range("AX197").Value = "Magnus Andersson"
$ws.Range("AY197").Value = "SCA Skog Naturvärdesinventering"

# Row 198
$ws.Range("A198").Value = 111973674
$ws.Range("B198").Value = 90652
$ws.Range("C198").Value = "Ovaliderad"
$ws.Range("D198").Value = "NT"
$ws.Range("E198").Value = 3100
$ws.Range("F198").Value = "Talltaggsvamp"
$ws.Range("G198").Value = "Bankera fuligineoalba"
$ws.Range("H198").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("I198").Value = "1"
$ws.Range("J198").Value = "mycel"
$ws.Range("P198").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q198").Value = 438160.5186564626
$ws.Range("R198").Value = 6953262.568950667
$ws.Range("S198").Value = 10
$ws.Range("T198").Value = "Jämtland"
$ws.Range("U198").Value = "Härjedalen"
$ws.Range("V198").Value = "Jämtland"
$ws.Range("W198").Value = "Vemdalen"
$ws.Range("Y198").Value = "2023-09-07"
$ws.Range("Z198").Value = "00:00"
$ws.Range("AA198").Value = "2023-09-07"
$ws.Range("AB198").Value = "00:00"
$ws.Range("AD198").Value = $false
$ws.Range("AE198").Value = $false
$ws.Range("AG198").Value = $false
$ws.Range("AI198").Value = "äldre renbetad fattigristallskog med lavfläck på torr moränmark"
$ws.Range("AN198").ClearContents()
$ws.Range("AO198").ClearContents()
$ws.Range("AT198").Value = ""
$ws.Range("AW198").Value = "Magnus Andersson"
$ws.Range("AX198").Value = "Magnus Andersson"
$ws.Range("AY198").Value = "SCA Skog Naturvärdesinventering"

# Row 199
$ws.Range("A199").Value = 111973673
$ws.Range("B199").Value = 90652
$ws.Range("C199").Value = "Ovaliderad"
$ws.Range("D199").Value = "NT"
$ws.Range("E199").Value = 3100
$ws.Range("F199").Value = "Talltaggsvamp"
$ws.Range("G199").Value = "Bankera fuligineoalba"
$ws.Range("H199").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("I199").Value = "1"
$ws.Range("J199").Value = "mycel"
$ws.Range("P199").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q199").Value = 438126.6098378488
$ws.Range("R199").Value = 6953243.924237223
$ws.Range("S199").Value = 10
$ws.Range("T199").Value = "Jämtland"
$ws.Range("U199").Value = "Härjedalen"
$ws.Range("V199").Value = "Jämtland"
$ws.Range("W199").Value = "Vemdalen"
$ws.Range("Y199").Value = "2023-09-07"
$ws.Range("Z199").Value = "00:00"
$ws.Range("AA199").Value = "2023-09-07"
$ws.Range("AB199").Value = "00:00"
$ws.Range("AD199").Value = $false
$ws.Range("AE199").Value = $false
$ws.Range("AG199").Value = $false
$ws.Range("AI199").Value = "äldre renbetad fattigristallskog med lavfläck på torr moränmark"
$ws.Range("AN199").ClearContents()
$ws.Range("AO199").ClearContents()
$ws.Range("AT199").Value = ""
$ws.Range("AW199").Value = "Magnus Andersson"
$ws.Range("AX199").Value = "Magnus Andersson"
$ws.Range("AY199").Value = "SCA Skog Naturvärdesinventering"

# Row 200
$ws.Range("A200").Value = 111973698
$ws.Range("B200").Value = 90660
$ws.Range("C200").Value = "Ovaliderad"
$ws.Range("D200").Value = "NT"
$ws.Range("E200").Value = 4362
$ws.Range("F200").Value = "Blå taggsvamp"
$ws.Range("G200").Value = "Hydnellum caeruleum"
$ws.Range("H200").Value = "(Hornem.) P.Karst."
$ws.Range("I200").Value = "1"
$ws.Range("J200").Value = "mycel"
$ws.Range("P200").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q200").Value = 438207.123851296
$ws.Range("R200").Value = 6953100.10165237
$ws.Range("S200").Value = 10
$ws.Range("T200").Value = "Jämtland"
$ws.Range("U200").Value = "Härjedalen"
$ws.Range("V200").Value = "Jämtland"
$ws.Range("W200").Value = "Vemdalen"
$ws.Range("Y200").Value = "2023-09-07"
$ws.Range("Z200").Value = "00:00"
$ws.Range("AA200").Value = "2023-09-07"
$ws.Range("AB200").Value = "00:00"
$ws.Range("AD200").Value = $false
$ws.Range("AE200").Value = $false
$ws.Range("AG200").Value = $false
$ws.Range("AI200").Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"
$ws.Range("AN200").ClearContents()
$ws.Range("AO200").ClearContents()
$ws.Range("AT200").Value = ""
$ws.Range("AW200").Value = "Magnus Andersson"
$ws.Range("AX200").Value = "Magnus Andersson"
$ws.Range("AY200").Value = "SCA Skog Naturvärdesinventering"

# Row 201
$ws.Range("A201").Value = 111973653
$ws.Range("B201").Value = 88032
$ws.Range("C201").Value = "Ovaliderad"
$ws.Range("D201").Value = "VU"
$ws.Range("E201").Value = 6276
$ws.Range("F201").Value = "Goliatmusseron"
$ws.Range("G201").Value = "Tricholoma matsutake"
$ws.Range("H201").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I201").Value = "1"
$ws.Range("J201").Value = "mycel"
$ws.Range("P201").Value = "Aloppmoarnas västra del, Jmt"
$ws.Range("Q201").Value = 438544.6107581231
$ws.Range("R201").Value = 6953546.421492521
$ws.Range("S201").Value = 10
$ws.Range("T201").Value = "Jämtland"
$ws.Range("U201").Value = "Härjedalen"
$ws.Range("V201").Value = "Jämtland"
$ws.Range("W201").Value = "Vemdalen"
$ws.Range("Y201").Value = "2023-09-07"
$ws.Range("Z201").Value = "00:00"
$ws.Range("AA201").Value = "2023-09-07"
$ws.Range("AB201").Value = "00:00"
$ws.Range("AD201").Value = $false
$ws.Range("AE201").Value = $false
$ws.Range("AG201").Value = $false
$ws.Range("AI201").Value = "äldre fattigristallskog på torr moränmark"
$ws.Range("AN201").ClearContents()
$ws.Range("AO201").ClearContents()
$ws.Range("AT201").Value = ""
$ws.Range("AW201").Value = "Magnus Andersson"
$ws.Range("AX201").Value = "Magnus Andersson"
$ws.Range("AY201").Value = "SCA Skog Naturvärdesinventering"
